$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("H5").Value = 3.3
$ws.Range("J5").Value = 1.08
$ws.Range("K5").Value = 8
$ws.Range("L5").Value = 1.44
$ws.Range("M5").Value = 2.63
$ws.Range("N5").Value = 2.35
$ws.Range("O5").Value = 1.57
$ws.Range("P5").Value = 1.53
$ws.Range("Q5").Value = 2.38
$ws.Range("R5").Value = 2
$ws.Range("S5").Value = 1.75
$ws.Range("T5").Value = 7.5
$ws.Range("Z5").Value = 8
$ws.Range("AC5").Value = 67
$ws.Range("AD5").Value = 451
$ws.Range("AE5").Value = 6.5
$ws.Range("AG5").Value = 10

# Row 7
$ws.Range("G7").Value = 4.2
$ws.Range("H7").Value = 3.75
$ws.Range("I7").Value = 1.8
$ws.Range("K7").Value = 13
$ws.Range("N7").Value = 1.85
$ws.Range("O7").Value = 2
$ws.Range("U7").Value = 21
$ws.Range("Y7").Value = 34
$ws.Range("AA7").Value = 7
$ws.Range("AC7").Value = 41
$ws.Range("AE7").Value = 8
$ws.Range("AF7").Value = 9
$ws.Range("AI7").Value = 15

# Row 9
$ws.Range("G9").Value = 1.24
$ws.Range("H9").Value = 5
$ws.Range("I9").Value = 11
$ws.Range("N9").Value = 1.85
$ws.Range("O9").Value = 2
$ws.Range("R9").Value = 2.38
$ws.Range("S9").Value = 1.53
$ws.Range("V9").Value = 9.5
$ws.Range("W9").Value = 7.5
$ws.Range("X9").Value = 12
$ws.Range("Y9").Value = 34
$ws.Range("AB9").Value = 26
$ws.Range("AE9").Value = 21
$ws.Range("AI9").Value = 81
$ws.Range("AJ9").Value = 81

# Row 17
$ws.Range("G17").Value = 2.75
$ws.Range("I17").Value = 2.63
$ws.Range("O17").Value = 1.36
$ws.Range("S17").Value = 1.47
$ws.Range("W17").Value = 29
$ws.Range("AH17").Value = 29
